$d = $word.ActiveDocument

# Collapse to the very end of the document body (just before sectPr's
# containing range end) so the new content lands after the last
# paragraph ("git rebase maser").
$r = $d.Content
$r.Collapse(0)

# Insert the new paragraph (with the smart/curly quotes and the
# gramStart/gramEnd proofing marks around "created") plus a trailing
# empty paragraph, using raw WordprocessingML so the run/proofErr
# structure matches exactly.
$w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$xml = '<w:p xmlns:w="' + $w + '">' +
       '<w:r><w:t xml:space="preserve">touch &#x201C;new file name to be </w:t></w:r>' +
       '<w:proofErr w:type="gramStart"/>' +
       '<w:r><w:t>created</w:t></w:r>' +
       '<w:proofErr w:type="gramEnd"/>' +
       '<w:r><w:t>&#x201D;</w:t></w:r>' +
       '</w:p>' +
       '<w:p xmlns:w="' + $w + '"/>'

$null = $r.InsertXML($xml)
